$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the last data row (row 61): thisletter -> "J", corrAns -> 1
$ws.Range("B61").Value = "J"
$ws.Range("C61").Value = 1

# Move the active selection to C61 (was E7)
$ws.Range("C61").Select()
